$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IndicatorMappingDB")
$ws.Activate()

# --- Update the "Tags" lookup values used by the Indicator Type ID / Name
#     columns (G:H) for the Interventions mapping rows.
#     These are shared strings, so updating the text once (via Replace)
#     propagates to every row that references them (rows 343-414).
$rng = $ws.Range("A1:R454")
$rng.Replace("1, 2, 3", "2, 3, 4", 1) | Out-Null
$rng.Replace("Drugs/Supplies, Labor, Visits", "Interventions, Labor, Visits", 1) | Out-Null

# --- Update the current selection on the IndicatorMappingDB sheet.
$ws.Range("P395").Select()
